$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 45.98144433333334
$ws.Range("H2").Value = 137.944333
$ws.Range("I2").Value = 0.9841234286873372
$ws.Range("J2").Value = 0.984123428687337
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.151158666666667
$ws.Range("N2").Value = 9.453476
$ws.Range("O2").Value = 0.03114707555614071
$ws.Range("P2").Value = 0.03114707555614071
$ws.Range("Q2").Value = 144.8948268168342
$ws.Range("R2").Value = 1304.053441351508
$ws.Range("S2").Value = 0.03065256678989275
$ws.Range("T2").Value = 0.03065256678989274

# Row 3
$ws.Range("G3").Value = 45.98144433333334
$ws.Range("H3").Value = 137.944333
$ws.Range("I3").Value = 0.9841234286873372
$ws.Range("J3").Value = 0.984123428687337
$ws.Range("M3").Value = 5.038243666666667
$ws.Range("O3").Value = 0.04979963650066307
$ws.Range("P3").Value = 0.04979963650066306
$ws.Range("Q3").Value = 231.6657206966026
$ws.Range("R3").Value = 2084.991486269423
$ws.Range("S3").Value = 0.0490089890204156
$ws.Range("T3").Value = 0.04900898902041559

# Row 4
$ws.Range("G4").Value = 45.98144433333334
$ws.Range("H4").Value = 137.944333
$ws.Range("I4").Value = 0.9841234286873372
$ws.Range("J4").Value = 0.984123428687337
$ws.Range("M4").Value = 92.911639
$ws.Range("N4").Value = 278.734917
$ws.Range("O4").Value = 0.9183688116343246
$ws.Range("P4").Value = 0.9183688116343246
$ws.Range("Q4").Value = 4272.211356597262
$ws.Range("R4").Value = 38449.90220937536
$ws.Range("S4").Value = 0.9037882637050868
$ws.Range("T4").Value = 0.9037882637050867

# Row 5
$ws.Range("G5").Value = 45.98144433333334
$ws.Range("H5").Value = 137.944333
$ws.Range("I5").Value = 0.9841234286873372
$ws.Range("J5").Value = 0.984123428687337
$ws.Range("M5").Value = 0.06924866666666667
$ws.Range("N5").Value = 0.207746
$ws.Range("O5").Value = 0.0006844763088715736
$ws.Range("P5").Value = 0.0006844763088715734
$ws.Range("Q5").Value = 3.184153711490889
$ws.Range("R5").Value = 28.657383403418
$ws.Range("S5").Value = 0.0006736091719419457
$ws.Range("T5").Value = 0.0006736091719419456

# Row 6
$ws.Range("I6").Value = 0.002244435796517234
$ws.Range("J6").Value = 0.002244435796517234
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 3.151158666666667
$ws.Range("N6").Value = 9.453476
$ws.Range("O6").Value = 0.03114707555614071
$ws.Range("P6").Value = 0.03114707555614071
$ws.Range("Q6").Value = 0.3304536062835556
$ws.Range("R6").Value = 2.974082456552
$ws.Range("S6").Value = 0.00006990761133502917
$ws.Range("T6").Value = 0.00006990761133502914

# Row 7
$ws.Range("I7").Value = 0.002244435796517234
$ws.Range("J7").Value = 0.002244435796517234
$ws.Range("M7").Value = 5.038243666666667
$ws.Range("O7").Value = 0.04979963650066307
$ws.Range("P7").Value = 0.04979963650066306
$ws.Range("Q7").Value = 0.5283471780068889
$ws.Range("R7").Value = 4.755124602062001
$ws.Range("S7").Value = 0.0001117720868156345
$ws.Range("T7").Value = 0.0001117720868156344

# Row 8
$ws.Range("I8").Value = 0.002244435796517234
$ws.Range("J8").Value = 0.002244435796517234
$ws.Range("M8").Value = 92.911639
$ws.Range("N8").Value = 278.734917
$ws.Range("O8").Value = 0.9183688116343246
$ws.Range("P8").Value = 0.9183688116343246
$ws.Range("Q8").Value = 9.743395817559332
$ws.Range("R8").Value = 87.690562358034
$ws.Range("S8").Value = 0.002061219835237071
$ws.Range("T8").Value = 0.002061219835237071

# Row 9
$ws.Range("I9").Value = 0.002244435796517234
$ws.Range("J9").Value = 0.002244435796517234
$ws.Range("M9").Value = 0.06924866666666667
$ws.Range("N9").Value = 0.207746
$ws.Range("O9").Value = 0.0006844763088715736
$ws.Range("P9").Value = 0.0006844763088715734
$ws.Range("Q9").Value = 0.007261923010222221
$ws.Range("R9").Value = 0.065357307092
$ws.Range("S9").Value = 0.000001536263129499347
$ws.Range("T9").Value = 0.000001536263129499346

# Row 10
$ws.Range("G10").Value = 0.547937
$ws.Range("H10").Value = 1.643811
$ws.Range("I10").Value = 0.01172728797372169
$ws.Range("J10").Value = 0.01172728797372169
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 3.151158666666667
$ws.Range("N10").Value = 9.453476
$ws.Range("O10").Value = 0.03114707555614071
$ws.Range("P10").Value = 0.03114707555614071
$ws.Range("Q10").Value = 1.726636426337333
$ws.Range("R10").Value = 15.539727837036
$ws.Range("S10").Value = 0.0003652707245861299
$ws.Range("T10").Value = 0.0003652707245861297

# Row 11
$ws.Range("G11").Value = 0.547937
$ws.Range("H11").Value = 1.643811
$ws.Range("I11").Value = 0.01172728797372169
$ws.Range("J11").Value = 0.01172728797372169
$ws.Range("M11").Value = 5.038243666666667
$ws.Range("O11").Value = 0.04979963650066307
$ws.Range("P11").Value = 0.04979963650066306
$ws.Range("Q11").Value = 2.760640119982334
$ws.Range("R11").Value = 24.845761079841
$ws.Range("S11").Value = 0.0005840146782299378
$ws.Range("T11").Value = 0.0005840146782299376

# Row 12
$ws.Range("G12").Value = 0.547937
$ws.Range("H12").Value = 1.643811
$ws.Range("I12").Value = 0.01172728797372169
$ws.Range("J12").Value = 0.01172728797372169
$ws.Range("M12").Value = 92.911639
$ws.Range("N12").Value = 278.734917
$ws.Range("O12").Value = 0.9183688116343246
$ws.Range("P12").Value = 0.9183688116343246
$ws.Range("Q12").Value = 50.90972473874299
$ws.Range("R12").Value = 458.187522648687
$ws.Range("S12").Value = 0.0107699755201203
$ws.Range("T12").Value = 0.0107699755201203

# Row 13
$ws.Range("G13").Value = 0.547937
$ws.Range("H13").Value = 1.643811
$ws.Range("I13").Value = 0.01172728797372169
$ws.Range("J13").Value = 0.01172728797372169
$ws.Range("M13").Value = 0.06924866666666667
$ws.Range("N13").Value = 0.207746
$ws.Range("O13").Value = 0.0006844763088715736
$ws.Range("P13").Value = 0.0006844763088715734
$ws.Range("Q13").Value = 0.03794390666733333
$ws.Range("R13").Value = 0.3414951600059999
$ws.Range("S13").Value = 0.00000802705078532702
$ws.Range("T13").Value = 0.000008027050785327015

# Row 14
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.08900066666666666
$ws.Range("H14").Value = 0.267002
$ws.Range("I14").Value = 0.001904847542424061
$ws.Range("J14").Value = 0.001904847542424061
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 3.151158666666667
$ws.Range("N14").Value = 9.453476
$ws.Range("O14").Value = 0.03114707555614071
$ws.Range("P14").Value = 0.03114707555614071
$ws.Range("Q14").Value = 0.2804552221057778
$ws.Range("R14").Value = 2.524096998952
$ws.Range("S14").Value = 0.0000593304303268112
$ws.Range("T14").Value = 0.00005933043032681118

# Row 15
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.08900066666666666
$ws.Range("H15").Value = 0.267002
$ws.Range("I15").Value = 0.001904847542424061
$ws.Range("J15").Value = 0.001904847542424061
$ws.Range("M15").Value = 5.038243666666667
$ws.Range("O15").Value = 0.04979963650066307
$ws.Range("P15").Value = 0.04979963650066306
$ws.Range("Q15").Value = 0.4484070451624445
$ws.Range("R15").Value = 4.035663406462001
$ws.Range("S15").Value = 0.00009486071520189964
$ws.Range("T15").Value = 0.00009486071520189961

# Row 16
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.08900066666666666
$ws.Range("H16").Value = 0.267002
$ws.Range("I16").Value = 0.001904847542424061
$ws.Range("J16").Value = 0.001904847542424061
$ws.Range("M16").Value = 92.911639
$ws.Range("N16").Value = 278.734917
$ws.Range("O16").Value = 0.9183688116343246
$ws.Range("P16").Value = 0.9183688116343246
$ws.Range("Q16").Value = 8.269197812092665
$ws.Range("R16").Value = 74.42278030883399
$ws.Range("S16").Value = 0.001749352573880549
$ws.Range("T16").Value = 0.001749352573880549

# Row 17
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.08900066666666666
$ws.Range("H17").Value = 0.267002
$ws.Range("I17").Value = 0.001904847542424061
$ws.Range("J17").Value = 0.001904847542424061
$ws.Range("M17").Value = 0.06924866666666667
$ws.Range("N17").Value = 0.207746
$ws.Range("O17").Value = 0.0006844763088715736
$ws.Range("P17").Value = 0.0006844763088715734
$ws.Range("Q17").Value = 0.006163177499111111
$ws.Range("R17").Value = 0.05546859749199999
$ws.Range("S17").Value = 0.00000130382301480151
$ws.Range("T17").Value = 0.000001303823014801509
